# Update from MV -datos- : revise last row (75) figures and append a new
# quarterly row (76) for 01-04-2021.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 75: revised figures (C75 is unchanged by this update) ---
$ws.Range("B75").Value = 15.5
$ws.Range("D75").Value = 1.2
$ws.Range("E75").Value = -138
$ws.Range("F75").Value = 119.6

# --- Row 76: new quarter ---
# Force the period label to be stored as text (matching the existing
# "dd-mm-yyyy"-look labels in column A) instead of being auto-converted
# to a date serial number, then drop the temporary number-format so the
# cell keeps the default (unstyled) look of the other data rows.
$lbl = $ws.Range("A76")
$lbl.NumberFormat = "@"
$lbl.Value = "01-04-2021"
$lbl.Style = "Normal"

$ws.Range("B76").Value = 15.4
$ws.Range("C76").Value = 15.3
$ws.Range("D76").Value = 0.6
$ws.Range("E76").Value = -127.8
$ws.Range("F76").Value = 113.3
